$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -4
$ws.Range("F3").Value = -6
$ws.Range("F5").Value = -5
$ws.Range("F6").Value = -2
$ws.Range("F8").Value = -5
$ws.Range("F11").Value = -2
$ws.Range("F12").Value = -6
$ws.Range("F14").Value = 2
$ws.Range("F23").Value = -8
$ws.Range("F25").Value = -4
$ws.Range("F26").Value = -2
$ws.Range("F27").Value = -4
$ws.Range("F33").Value = -4
